# Process migration to LT-00179
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update refreshed Market Cap values (shared string content changes)
$ws.Range("C2").Value = "₹ 308,461 Cr."   # WIPRO
$ws.Range("C3").Value = "₹ 38,372 Cr."    # IOB
$ws.Range("C5").Value = "₹ 552,371 Cr."   # ICICI
$ws.Range("C6").Value = "₹ 173,189 Cr."   # TATA
$ws.Range("C8").Value = "₹ 207,386 Cr."   # ONGC
$ws.Range("C9").Value = "₹ 38,015 Cr."    # ASHOKLEY

# Replace ASIANPAINT entry with SLICE, which now has a company code discrepancy
$ws.Range("A11").Value = "SLICE"
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "-"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "Company code discrepancy identified"
